$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap full row content (C:L) between row 12 and row 13 ---
# (cams-global-radiative-forcing-auxilliary-variables <-> cams-global-radiative-forcings)
foreach ($col in @("C","D","E","F","G","H","I","J","K","L")) {
    $addr12 = $col + "12"
    $addr13 = $col + "13"
    $v12 = $ws.Range($addr12).Value2
    $v13 = $ws.Range($addr13).Value2
    $ws.Range($addr12).Value = $v13
    $ws.Range($addr13).Value = $v12
}

# --- Swap full row content (C:I, K:L) between row 15 and row 17 ---
# (cams-solar-radiation-timeseries <-> cams-global-reanalysis-eac4-monthly); J handled separately below
foreach ($col in @("C","D","E","F","G","H","I","K","L")) {
    $addr15 = $col + "15"
    $addr17 = $col + "17"
    $v15 = $ws.Range($addr15).Value2
    $v17 = $ws.Range($addr17).Value2
    $ws.Range($addr15).Value = $v17
    $ws.Range($addr17).Value = $v15
}

# --- Explicit "extent" (column J) date refresh for all affected rows ---
$ws.Range("J2").Value = "{'spatial': {'bbox': [[0, -89, 360, 89]]}, 'temporal': {'interval': [['2024-03-01T00:00:00Z', '2026-02-07T00:00:00Z']]}}"
$ws.Range("J6").Value = "{'spatial': {'bbox': [[-25, 30, 45, 72]]}, 'temporal': {'interval': [['2024-01-17T00:00:00Z', '2026-02-07T00:00:00Z']]}}"
$ws.Range("J11").Value = "{'spatial': {'bbox': [[-180, -90, 180, 90]]}, 'temporal': {'interval': [['2015-01-01T00:00:00Z', '2026-02-07T00:00:00Z']]}}"
$ws.Range("J15").Value = "{'spatial': {'bbox': [[-180, -90, 180, 90]]}, 'temporal': {'interval': [['2003-01-01T00:00:00Z', '2023-12-31T00:00:00Z']]}}"
$ws.Range("J16").Value = "{'spatial': {'bbox': [[-25, 30, 45, 72]]}, 'temporal': {'interval': [['2023-02-02T00:00:00Z', '2026-02-08T00:00:00Z']]}}"
$ws.Range("J17").Value = "{'spatial': {'bbox': [[0, -89, 360, 89]]}, 'temporal': {'interval': [['2004-01-01T00:00:00Z', '2026-02-06T00:00:00Z']]}}"
